# Update Metadata sheet values
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.4.0-snapshot-1"
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# Swap the two "Mapping" columns (AK <-> AL) on the Elements sheet, including
# header text, data, and column widths. (AK was the narrow "RIM Mapping"
# column, AL was the wide "Spécification métier" column; after the edit
# they trade places.)
$elements = $wb.Worksheets.Item("Elements")

for ($row = 1; $row -le 6; $row++) {
    $akCell = $elements.Cells.Item($row, 37)
    $alCell = $elements.Cells.Item($row, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Column AK (37) takes AL's former (wide) width; AL (38) takes AK's former
# (narrow) width.
$elements.Columns(37).ColumnWidth = 84.83
$elements.Columns(38).ColumnWidth = 24.16
